$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "67.126.66"
Set-TextValue $ws.Range("E2") "  +3.53%  "
Set-TextValue $ws.Range("D3") "3.856.91"
Set-TextValue $ws.Range("E3") "  +5.61%  "
Set-TextValue $ws.Range("E4") "  -0.23%  "
Set-TextValue $ws.Range("D5") "423.73"
Set-TextValue $ws.Range("E5") "  +4.12%  "
Set-TextValue $ws.Range("D6") "129.51"
Set-TextValue $ws.Range("E6") "  -3.01%  "
Set-TextValue $ws.Range("D7") "3.851.96"
Set-TextValue $ws.Range("E7") "  +5.73%  "
Set-TextValue $ws.Range("D8") "0.609"
Set-TextValue $ws.Range("E8") "  -1.78%  "
Set-TextValue $ws.Range("D9") "0.999"
Set-TextValue $ws.Range("E9") "  -0.17%  "
Set-TextValue $ws.Range("D10") "0.724"
Set-TextValue $ws.Range("E10") "  -0.52%  "
Set-TextValue $ws.Range("D11") "0.158"
Set-TextValue $ws.Range("E11") "  -2.47%  "
Set-TextValue $ws.Range("D12") "0.0000342"
Set-TextValue $ws.Range("E12") "  +6.63%  "
Set-TextValue $ws.Range("D13") "40.90"
Set-TextValue $ws.Range("E13") "  -2.95%  "
Set-TextValue $ws.Range("B14") "Polkadot"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D14") "10.40"
Set-TextValue $ws.Range("E14") "  +4.74%  "
Set-TextValue $ws.Range("B15") "WrappedliquidstakedEther2.0"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D15") "4.461.01"
Set-TextValue $ws.Range("E15") "  +4.81%  "
Set-TextValue $ws.Range("D16") "15.72"
Set-TextValue $ws.Range("E16") "  +15.12%  "
Set-TextValue $ws.Range("D17") "3.854.61"
Set-TextValue $ws.Range("E17") "  +5.49%  "
Set-TextValue $ws.Range("E18") "  -0.59%  "
Set-TextValue $ws.Range("D19") "19.84"
Set-TextValue $ws.Range("E19") "  -0.85%  "
Set-TextValue $ws.Range("D20") "67.360.16"
Set-TextValue $ws.Range("E20") "  +3.64%  "
Set-TextValue $ws.Range("E21") "  -0.17%  "
Set-TextValue $ws.Range("D22") "409.57"
Set-TextValue $ws.Range("E22") "  -2.95%  "
Set-TextValue $ws.Range("D23") "15.00"
Set-TextValue $ws.Range("E23") "  -1.80%  "
Set-TextValue $ws.Range("D24") "84.27"
Set-TextValue $ws.Range("E24") "  -1.97%  "
Set-TextValue $ws.Range("E25") "  +1.50%  "
Set-TextValue $ws.Range("D26") "37.41"
Set-TextValue $ws.Range("E26") "  +4.44%  "
Set-TextValue $ws.Range("D27") "10.02"
Set-TextValue $ws.Range("E27") "  +6.60%  "
Set-TextValue $ws.Range("D28") "3.24"
Set-TextValue $ws.Range("E28") "  +1.06%  "
Set-TextValue $ws.Range("B29") "RenderToken"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D29") "9.53"
Set-TextValue $ws.Range("E29") "  +36.89%  "
Set-TextValue $ws.Range("B30") "LEO"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D30") "5.43"
Set-TextValue $ws.Range("E30") "  +5.87%  "
Set-TextValue $ws.Range("D31") "736.72"
Set-TextValue $ws.Range("E31") "  +8.47%  "
Set-TextValue $ws.Range("D32") "13.16"
Set-TextValue $ws.Range("E32") "  +3.46%  "
Set-TextValue $ws.Range("E33") "  +2.14%  "
Set-TextValue $ws.Range("E34") "  +3.23%  "
Set-TextValue $ws.Range("E35") "  -0.06%  "
Set-TextValue $ws.Range("E36") "  -5.98%  "
Set-TextValue $ws.Range("D37") "38.60"
Set-TextValue $ws.Range("E37") "  -7.17%  "
Set-TextValue $ws.Range("D38") "55.67"
Set-TextValue $ws.Range("E38") "  -0.55%  "
Set-TextValue $ws.Range("E39") "  +22.37%  "
Set-TextValue $ws.Range("D40") "0.0₃0748"
Set-TextValue $ws.Range("E40") "  +15.14%  "
Set-TextValue $ws.Range("D41") "0.0457"
Set-TextValue $ws.Range("E41") "  -1.79%  "
Set-TextValue $ws.Range("D42") "2.90"
Set-TextValue $ws.Range("E42") "  -1.99%  "
Set-TextValue $ws.Range("E43") "  +0.46%  "
Set-TextValue $ws.Range("D44") "3.37"
Set-TextValue $ws.Range("E44") "  +0.92%  "
Set-TextValue $ws.Range("E45") "  -4.31%  "
Set-TextValue $ws.Range("B46") "TheGraph"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D46") "0.316"
Set-TextValue $ws.Range("E46") "  +8.70%  "
Set-TextValue $ws.Range("B47") "ApeXProtocol"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws.Range("D47") "3.11"
Set-TextValue $ws.Range("E47") "  -0.25%  "
Set-TextValue $ws.Range("B48") "ARBITRUM"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D48") "2.05"
Set-TextValue $ws.Range("E48") "  -1.58%  "
Set-TextValue $ws.Range("B49") "Monero"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D49") "141.91"
Set-TextValue $ws.Range("E49") "  -1.32%  "
Set-TextValue $ws.Range("D50") "2.83"
Set-TextValue $ws.Range("E50") "  +0.30%  "
Set-TextValue $ws.Range("B51") "EnergySwap"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D51") "25.67"
Set-TextValue $ws.Range("E51") "  -3.35%  "
